$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.215.71"
$ws.Range("E2").Value = "  +1.98%  "
$ws.Range("D3").Value = "3.322.18"
$ws.Range("E3").Value = "  +2.70%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.16"
$ws.Range("E5").Value = "  +3.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.57"
$ws.Range("E6").Value = "  +3.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "3.316.07"
$ws.Range("E8").Value = "  +2.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.577"
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.177"
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.572"
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.20"
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000265"
$ws.Range("E13").Value = "  +2.62%  "
$ws.Range("D14").Value = "3.845.88"
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.47"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "612.24"
$ws.Range("E16").Value = "  -1.54%  "
$ws.Range("D17").Value = "66.201.84"
$ws.Range("E17").Value = "  +1.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.118"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.82"
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("B20").Value = "WrappedEther"
$ws.Range("C20").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").Value = "3.321.42"
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.00"
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.899"
$ws.Range("E22").Value = "  +1.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.21"
$ws.Range("E23").Value = "  +3.66%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "100.09"
$ws.Range("E24").Value = "  -4.95%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.00"
$ws.Range("E25").Value = "  +4.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.00"
$ws.Range("E26").Value = "  +4.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.72"
$ws.Range("E27").Value = "  +4.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.48"
$ws.Range("E28").Value = "  +2.51%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "30.99"
$ws.Range("E29").Value = "  +5.67%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.51"
$ws.Range("E30").Value = "  +1.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.60"
$ws.Range("E31").Value = "  +8.44%  "
$ws.Range("B32").Value = "dogwifhat"
$ws.Range("C32").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.71"
$ws.Range("E32").Value = "  -0.62%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "561.99"
$ws.Range("E33").Value = "  +7.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.89"
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("D35").Value = "3.830.56"
$ws.Range("E35").Value = "  +3.11%  "
$ws.Range("E36").Value = "  +1.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "56.12"
$ws.Range("E38").Value = "  -0.82%  "
$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "33.34"
$ws.Range("E39").Value = "  +4.29%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.128"
$ws.Range("E40").Value = "  +1.35%  "
$ws.Range("B41").Value = "ApeXProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.42"
$ws.Range("E41").Value = "  +6.22%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.15"
$ws.Range("E42").Value = "  -1.52%  "
$ws.Range("D43").Value = "0.0₃0686"
$ws.Range("E43").Value = "  -2.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.60"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.332"
$ws.Range("E45").Value = "  +2.31%  "
$ws.Range("E46").Value = "  +1.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.09"
$ws.Range("E47").Value = "  -3.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.127"
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.54"
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "128.93"
$ws.Range("E51").Value = "  +6.16%  "
